$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PE, JANUARY")

# Rows 15-18: set the purchase date text (01/30/2026) while preserving the
# existing date-formatted cell style (numFmt 14 / bordered), the same way
# row 14 (01/29/2026) is already stored - as a literal shared string, not
# a parsed date serial. We do this by entering the text under a temporary
# "Text" number format, then pasting back the original cell's format
# (copied from C14) without touching the value.
foreach ($r in 15..18) {
    $dst = $ws.Range("C$r")
    $dst.NumberFormat = "@"
    $dst.Value = "01/30/2026"
    $ws.Range("C14").Copy()
    $dst.PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Invoice / reference numbers (column G)
$ws.Range("G15").Value = 518341441
$ws.Range("G16").Value = 518341067
$ws.Range("G17").Formula = "=518341060"
$ws.Range("G18").Value = 518341429

# Gross purchases (column I) - net of input tax, entered as formulas
$ws.Range("I15").Formula = "=1232820-49636.85"
$ws.Range("I16").Formula = "=154728-19963.74"
$ws.Range("I17").Formula = "=794772-176047.41"
$ws.Range("I18").Formula = "=44822-3731.76"

$ws.Range("C19").Select()

$wb.Save()
